$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the payer name in B7 from "WINNIE SAVIRA" to "WINNIE SAVIRA1"
$ws.Range("B7").Value = "WINNIE SAVIRA1"

# Update the active cell / selection to B10
$ws.Range("B10").Select()
